# Adds condition-specificity correlation columns (AR:AV) to the cluster
# enrichment table: corr_rest, corr_stim8hr, corr_stim48hr, corr_shared,
# and a condition_specificity category label per cluster row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AR:AV ---
$ws.Range("AR1").Value = "corr_rest"
$ws.Range("AS1").Value = "corr_stim8hr"
$ws.Range("AT1").Value = "corr_stim48hr"
$ws.Range("AU1").Value = "corr_shared"
$ws.Range("AV1").Value = "condition_specificity"

# --- Data rows 2-113 ---
# AR = corr_rest, AS = corr_stim8hr, AT = corr_stim48hr, AU = corr_shared,
# AV = condition_specificity (across_condition / Rest / Stim8hr / Stim48hr / ...)
$ws.Range("AR2").Value = [double]"0.40648157264866502"
$ws.Range("AS2").Value = [double]"0.51297608215929202"
$ws.Range("AT2").Value = [double]"0.49674284784620598"
$ws.Range("AU2").Value = [double]"0.35607590946914303"
$ws.Range("AV2").Value = "across_condition"
$ws.Range("AR3").Value = [double]"0.38271143951381198"
$ws.Range("AS3").Value = [double]"0.41156398933175597"
$ws.Range("AT3").Value = [double]"0.29206693582232202"
$ws.Range("AU3").Value = [double]"0.327720702663717"
$ws.Range("AV3").Value = "across_condition"
$ws.Range("AR4").Value = [double]"0.26168032925378898"
$ws.Range("AS4").Value = [double]"0.304253618567906"
$ws.Range("AT4").Value = [double]"0.28165853621510001"
$ws.Range("AU4").Value = [double]"0.20715558396926001"
$ws.Range("AV4").Value = "across_condition"
$ws.Range("AR5").Value = [double]"0.50855539183175003"
$ws.Range("AS5").Value = [double]"0.50356974752498596"
$ws.Range("AT5").Value = [double]"0.45303808566172399"
$ws.Range("AU5").Value = [double]"0.29677971906460099"
$ws.Range("AV5").Value = "across_condition"
$ws.Range("AR6").Value = [double]"0.49469050661269798"
$ws.Range("AS6").Value = [double]"0.51875399674508804"
$ws.Range("AT6").Value = [double]"0.53517958449020597"
$ws.Range("AU6").Value = [double]"0.37651358885085501"
$ws.Range("AV6").Value = "across_condition"
$ws.Range("AR7").Value = [double]"0.58384714986211395"
$ws.Range("AS7").Value = [double]"0.59607816794851798"
$ws.Range("AT7").Value = [double]"0.59594003507703197"
$ws.Range("AU7").Value = [double]"0.45855450764099598"
$ws.Range("AV7").Value = "across_condition"
$ws.Range("AR8").Value = [double]"0.28189776494507501"
$ws.Range("AS8").Value = [double]"0.28429985881269598"
$ws.Range("AT8").Value = [double]"0.20225040592572199"
$ws.Range("AU8").Value = [double]"0.16879182340447499"
$ws.Range("AV8").Value = "across_condition"
$ws.Range("AR9").Value = [double]"0.37561466525575599"
$ws.Range("AS9").Value = [double]"0.41273355039861997"
$ws.Range("AT9").Value = [double]"0.25022905039843402"
$ws.Range("AU9").Value = [double]"0.27357749524192898"
$ws.Range("AV9").Value = "across_condition"
$ws.Range("AR10").Value = [double]"0.287589695131717"
$ws.Range("AS10").Value = [double]"0.34453305238339998"
$ws.Range("AT10").Value = [double]"0.27203185682910003"
$ws.Range("AU10").Value = [double]"0.156596972755804"
$ws.Range("AV10").Value = "across_condition"
$ws.Range("AR11").Value = [double]"0.35403665293375902"
$ws.Range("AS11").Value = [double]"0.43224428218294902"
$ws.Range("AT11").Value = [double]"0.37565231970179402"
$ws.Range("AU11").Value = [double]"0.245895864279606"
$ws.Range("AV11").Value = "across_condition"
$ws.Range("AR12").Value = [double]"0.311337690331072"
$ws.Range("AS12").Value = [double]"0.40465009802145702"
$ws.Range("AT12").Value = [double]"0.29021946381957803"
$ws.Range("AU12").Value = [double]"0.20646226366123299"
$ws.Range("AV12").Value = "across_condition"
$ws.Range("AR13").Value = [double]"0.75538941529859105"
$ws.Range("AS13").Value = [double]"0.79334963943351999"
$ws.Range("AT13").Value = [double]"0.71902620140918105"
$ws.Range("AU13").Value = [double]"0.451719145779926"
$ws.Range("AV13").Value = "across_condition"
$ws.Range("AR14").Value = [double]"0.24250318930688"
$ws.Range("AS14").Value = [double]"0.34385338412932498"
$ws.Range("AT14").Value = [double]"0.35590591765028901"
$ws.Range("AU14").Value = [double]"0.21874806608835301"
$ws.Range("AV14").Value = "across_condition"
$ws.Range("AR15").Value = [double]"0.48218023965294599"
$ws.Range("AS15").Value = [double]"0.49235093336580299"
$ws.Range("AT15").Value = [double]"0.34789354864596"
$ws.Range("AU15").Value = [double]"0.30127929243325402"
$ws.Range("AV15").Value = "across_condition"
$ws.Range("AR16").Value = [double]"0.268236380429831"
$ws.Range("AS16").Value = [double]"0.357800520734325"
$ws.Range("AT16").Value = [double]"0.19875531167617599"
$ws.Range("AU16").Value = [double]"0.12201946751661499"
$ws.Range("AV16").Value = "across_condition"
$ws.Range("AR17").Value = [double]"0.25503408921506399"
$ws.Range("AS17").Value = [double]"0.37283674663944599"
$ws.Range("AT17").Value = [double]"0.37844895508888199"
$ws.Range("AU17").Value = [double]"0.13771904465812199"
$ws.Range("AV17").Value = "across_condition"
$ws.Range("AR18").Value = [double]"0.36979401823032998"
$ws.Range("AS18").Value = [double]"0.38330684704201501"
$ws.Range("AT18").Value = [double]"0.38416212936455901"
$ws.Range("AU18").Value = [double]"0.271141655204794"
$ws.Range("AV18").Value = "across_condition"
$ws.Range("AR19").Value = [double]"0.25528761782811898"
$ws.Range("AS19").Value = [double]"0.23844987720317901"
$ws.Range("AT19").Value = [double]"0.26410021817524099"
$ws.Range("AU19").Value = [double]"0.207045604286653"
$ws.Range("AV19").Value = "across_condition"
$ws.Range("AR20").Value = [double]"0.39235980654043201"
$ws.Range("AS20").Value = [double]"0.49322054874429799"
$ws.Range("AT20").Value = [double]"0.55545138138686601"
$ws.Range("AU20").Value = [double]"0.36251967020011999"
$ws.Range("AV20").Value = "across_condition"
$ws.Range("AR21").Value = [double]"0.20741186107882001"
$ws.Range("AS21").Value = [double]"0.199395605586102"
$ws.Range("AT21").Value = [double]"0.18333538407161301"
$ws.Range("AU21").Value = [double]"0.115186610110244"
$ws.Range("AV21").Value = "across_condition"
$ws.Range("AR22").Value = [double]"0.36252083442198801"
$ws.Range("AS22").Value = [double]"0.27371662107758099"
$ws.Range("AT22").Value = [double]"0.32442984803907698"
$ws.Range("AU22").Value = [double]"0.23232120711534299"
$ws.Range("AV22").Value = "across_condition"
$ws.Range("AR23").Value = [double]"3.3504841932685402E-2"
$ws.Range("AS23").Value = [double]"0.25087596509236298"
$ws.Range("AT23").Value = [double]"0.113018943554052"
$ws.Range("AU23").Value = [double]"7.13760208017446E-2"
$ws.Range("AV23").Value = "Stim8hr"
$ws.Range("AR24").Value = [double]"0.172067814723216"
$ws.Range("AS24").Value = [double]"0.900215029548178"
$ws.Range("AT24").Value = [double]"0.62014296741337305"
$ws.Range("AU24").Value = [double]"0.29636035489868601"
$ws.Range("AV24").Value = "Stim8hr_Stim48hr"
$ws.Range("AR25").Value = [double]"0.28849869005652301"
$ws.Range("AS25").Value = [double]"0.34349970376639899"
$ws.Range("AT25").Value = [double]"0.20409895873579101"
$ws.Range("AU25").Value = [double]"0.17415918154028001"
$ws.Range("AV25").Value = "across_condition"
$ws.Range("AR26").Value = [double]"0.110588429405108"
$ws.Range("AS26").Value = [double]"0.22876910275840201"
$ws.Range("AT26").Value = [double]"0.11544253539788001"
$ws.Range("AU26").Value = [double]"0.116910295167101"
$ws.Range("AV26").Value = "across_condition"
$ws.Range("AR27").Value = [double]"0.18042577592551301"
$ws.Range("AS27").Value = [double]"7.7102053545282795E-2"
$ws.Range("AT27").Value = [double]"8.0545005470863601E-2"
$ws.Range("AU27").Value = [double]"4.6398651408720799E-2"
$ws.Range("AV27").Value = "Rest"
$ws.Range("AR28").Value = [double]"0.18350073973968201"
$ws.Range("AS28").Value = [double]"3.4000560144319597E-2"
$ws.Range("AT28").Value = [double]"2.4623304586057299E-2"
$ws.Range("AU28").Value = [double]"2.9672100746200601E-2"
$ws.Range("AV28").Value = "Rest"
$ws.Range("AR29").Value = [double]"0.38746167346253202"
$ws.Range("AS29").Value = [double]"8.7457937143009506E-2"
$ws.Range("AT29").Value = [double]"4.3984120189526701E-2"
$ws.Range("AU29").Value = [double]"6.89300535638252E-2"
$ws.Range("AV29").Value = "Rest"
$ws.Range("AR30").Value = [double]"0.29526807576053699"
$ws.Range("AS30").Value = [double]"0.37806008357562898"
$ws.Range("AT30").Value = [double]"0.41565626936480299"
$ws.Range("AU30").Value = [double]"0.15783483621658101"
$ws.Range("AV30").Value = "across_condition"
$ws.Range("AR31").Value = [double]"0.243426488047329"
$ws.Range("AS31").Value = [double]"0.10542800740903099"
$ws.Range("AT31").Value = [double]"8.7875972386729798E-2"
$ws.Range("AU31").Value = [double]"9.1543584829030702E-2"
$ws.Range("AV31").Value = "Rest"
$ws.Range("AR32").Value = [double]"0.224973509299247"
$ws.Range("AS32").Value = [double]"0.19122984520759201"
$ws.Range("AT32").Value = [double]"6.63446240441697E-2"
$ws.Range("AU32").Value = [double]"4.77887677872013E-2"
$ws.Range("AV32").Value = "Rest"
$ws.Range("AR33").Value = [double]"0.250172878152711"
$ws.Range("AS33").Value = [double]"0.30362092027127802"
$ws.Range("AT33").Value = [double]"0.41401639969723297"
$ws.Range("AU33").Value = [double]"0.17752080633779399"
$ws.Range("AV33").Value = "across_condition"
$ws.Range("AR34").Value = [double]"0.251501747177026"
$ws.Range("AS34").Value = [double]"0.36115155966906798"
$ws.Range("AT34").Value = [double]"0.372519905534624"
$ws.Range("AU34").Value = [double]"0.191224909942628"
$ws.Range("AV34").Value = "across_condition"
$ws.Range("AR35").Value = [double]"0.39622177618814403"
$ws.Range("AS35").Value = [double]"0.42842397973617302"
$ws.Range("AT35").Value = [double]"0.27355318881551799"
$ws.Range("AU35").Value = [double]"0.21671154919640601"
$ws.Range("AV35").Value = "across_condition"
$ws.Range("AR36").Value = [double]"0.28437150705320302"
$ws.Range("AS36").Value = [double]"0.35165615023746599"
$ws.Range("AT36").Value = [double]"0.27791562295181799"
$ws.Range("AU36").Value = [double]"0.16138126127505301"
$ws.Range("AV36").Value = "across_condition"
$ws.Range("AR37").Value = [double]"0.20697038339519999"
$ws.Range("AS37").Value = [double]"0.215302123842608"
$ws.Range("AT37").Value = [double]"0.15472007448007699"
$ws.Range("AU37").Value = [double]"9.7159515932049204E-2"
$ws.Range("AV37").Value = "across_condition"
$ws.Range("AR38").Value = [double]"7.5987417075073102E-2"
$ws.Range("AS38").Value = [double]"8.6317682594838796E-2"
$ws.Range("AT38").Value = [double]"0.39753842683851598"
$ws.Range("AU38").Value = [double]"6.4638964075712096E-2"
$ws.Range("AV38").Value = "Stim48hr"
$ws.Range("AR39").Value = [double]"0.20352167831267001"
$ws.Range("AS39").Value = [double]"0.29061855031836298"
$ws.Range("AT39").Value = [double]"7.7771240142422704E-2"
$ws.Range("AU39").Value = [double]"0.115547529397677"
$ws.Range("AV39").Value = "Rest_Stim8hr"
$ws.Range("AR40").Value = [double]"0.23122864775617699"
$ws.Range("AS40").Value = [double]"0.267279913130827"
$ws.Range("AT40").Value = [double]"0.32149989690660102"
$ws.Range("AU40").Value = [double]"0.16949293678261601"
$ws.Range("AV40").Value = "across_condition"
$ws.Range("AR41").Value = [double]"0.134536090906998"
$ws.Range("AS41").Value = [double]"0.14010723758882401"
$ws.Range("AT41").Value = [double]"0.13399534147001699"
$ws.Range("AU41").Value = [double]"9.0479769947255401E-2"
$ws.Range("AV41").Value = "across_condition"
$ws.Range("AR42").Value = [double]"0.186772062054714"
$ws.Range("AS42").Value = [double]"0.175079787843753"
$ws.Range("AT42").Value = [double]"0.25768093288912097"
$ws.Range("AU42").Value = [double]"0.12836132264198299"
$ws.Range("AV42").Value = "across_condition"
$ws.Range("AR43").Value = [double]"0.16152240835103901"
$ws.Range("AS43").Value = [double]"0.20250903713642901"
$ws.Range("AT43").Value = [double]"8.2737605975873804E-2"
$ws.Range("AU43").Value = [double]"0.11389618134115"
$ws.Range("AV43").Value = "across_condition"
$ws.Range("AR44").Value = [double]"0.16674749886375401"
$ws.Range("AS44").Value = [double]"0.13171643414008"
$ws.Range("AT44").Value = [double]"5.6037628746958203E-2"
$ws.Range("AU44").Value = [double]"6.8663211981368794E-2"
$ws.Range("AV44").Value = "across_condition"
$ws.Range("AR45").Value = [double]"0.34845211077336302"
$ws.Range("AS45").Value = [double]"0.37356738862746602"
$ws.Range("AT45").Value = [double]"0.41849434646970501"
$ws.Range("AU45").Value = [double]"0.228452239957908"
$ws.Range("AV45").Value = "across_condition"
$ws.Range("AR46").Value = [double]"0.22972189862016201"
$ws.Range("AS46").Value = [double]"0.27931444905078001"
$ws.Range("AT46").Value = [double]"0.29921232871306602"
$ws.Range("AU46").Value = [double]"0.16489458507812399"
$ws.Range("AV46").Value = "across_condition"
$ws.Range("AR47").Value = [double]"7.8767443366201104E-2"
$ws.Range("AS47").Value = [double]"3.6334667614455897E-2"
$ws.Range("AT47").Value = [double]"3.2501420696840901E-2"
$ws.Range("AU47").Value = [double]"3.8355007838106699E-2"
$ws.Range("AV47").Value = "Rest"
$ws.Range("AR48").Value = [double]"0.38543682840562299"
$ws.Range("AS48").Value = [double]"0.36959149065069402"
$ws.Range("AT48").Value = [double]"0.24127164564936701"
$ws.Range("AU48").Value = [double]"0.13593278050574401"
$ws.Range("AV48").Value = "across_condition"
$ws.Range("AR49").Value = [double]"0.373931850626191"
$ws.Range("AS49").Value = [double]"0.33861540451188898"
$ws.Range("AT49").Value = [double]"0.32563880836249298"
$ws.Range("AU49").Value = [double]"0.257956346514153"
$ws.Range("AV49").Value = "across_condition"
$ws.Range("AR50").Value = [double]"0.32968273420414701"
$ws.Range("AS50").Value = [double]"0.22887767139457299"
$ws.Range("AT50").Value = [double]"0.30003147425503002"
$ws.Range("AU50").Value = [double]"0.15552176616518101"
$ws.Range("AV50").Value = "across_condition"
$ws.Range("AR51").Value = [double]"0.215605048700966"
$ws.Range("AS51").Value = [double]"3.9337523424069001E-2"
$ws.Range("AT51").Value = [double]"1.25448841308342E-2"
$ws.Range("AU51").Value = [double]"2.83840268812376E-2"
$ws.Range("AV51").Value = "Rest"
$ws.Range("AR52").Value = [double]"0.16156134884538101"
$ws.Range("AS52").Value = [double]"0.148261743038378"
$ws.Range("AT52").Value = [double]"0.161162129693531"
$ws.Range("AU52").Value = [double]"5.7857263118178598E-2"
$ws.Range("AV52").Value = "across_condition"
$ws.Range("AR53").Value = [double]"0.21476975991448699"
$ws.Range("AS53").Value = [double]"0.116426118066072"
$ws.Range("AT53").Value = [double]"0.14942735973636001"
$ws.Range("AU53").Value = [double]"5.6503926523180897E-2"
$ws.Range("AV53").Value = "across_condition"
$ws.Range("AR54").Value = [double]"0.34297196646785399"
$ws.Range("AS54").Value = [double]"0.320333345215705"
$ws.Range("AT54").Value = [double]"0.169974847864198"
$ws.Range("AU54").Value = [double]"0.13967692641989399"
$ws.Range("AV54").Value = "across_condition"
$ws.Range("AR55").Value = [double]"9.1972968582308301E-2"
$ws.Range("AS55").Value = [double]"0.187582977509199"
$ws.Range("AT55").Value = [double]"6.3077173217919003E-2"
$ws.Range("AU55").Value = [double]"5.3531723514045899E-2"
$ws.Range("AV55").Value = "Stim8hr"
$ws.Range("AR56").Value = [double]"0.25139559160368102"
$ws.Range("AS56").Value = [double]"0.23161495180429401"
$ws.Range("AT56").Value = [double]"5.0932947357760001E-2"
$ws.Range("AU56").Value = [double]"6.2277401798033601E-2"
$ws.Range("AV56").Value = "Rest_Stim8hr"
$ws.Range("AR57").Value = [double]"-4.9630384116220003E-4"
$ws.Range("AS57").Value = [double]"-3.8964381905269498E-2"
$ws.Range("AT57").Value = [double]"1.9475633263060201E-2"
$ws.Range("AU57").Value = [double]"6.983751474702E-4"
$ws.Range("AV57").Value = "across_condition"
$ws.Range("AR58").Value = [double]"0.62813684169987405"
$ws.Range("AS58").Value = [double]"0.50457297407731105"
$ws.Range("AT58").Value = [double]"0.39741083095810298"
$ws.Range("AU58").Value = [double]"0.364886345802932"
$ws.Range("AV58").Value = "across_condition"
$ws.Range("AR59").Value = [double]"7.92065395484782E-2"
$ws.Range("AS59").Value = [double]"6.2758499077573193E-2"
$ws.Range("AT59").Value = [double]"0.31049389368074998"
$ws.Range("AU59").Value = [double]"5.0583446267173898E-2"
$ws.Range("AV59").Value = "Stim48hr"
$ws.Range("AR60").Value = [double]"0.18328800279784399"
$ws.Range("AS60").Value = [double]"0.18756449818594401"
$ws.Range("AT60").Value = [double]"0.37904466656708602"
$ws.Range("AU60").Value = [double]"9.8190380542557701E-2"
$ws.Range("AV60").Value = "across_condition"
$ws.Range("AR61").Value = [double]"0.121364227429205"
$ws.Range("AS61").Value = [double]"0.101895296806969"
$ws.Range("AT61").Value = [double]"0.121308187677632"
$ws.Range("AU61").Value = [double]"7.1749672156205299E-2"
$ws.Range("AV61").Value = "across_condition"
$ws.Range("AR62").Value = [double]"0.24522877116095901"
$ws.Range("AS62").Value = [double]"0.35674486540871603"
$ws.Range("AT62").Value = [double]"9.5396082361758394E-2"
$ws.Range("AU62").Value = [double]"6.0240799549630297E-2"
$ws.Range("AV62").Value = "Rest_Stim8hr"
$ws.Range("AR63").Value = [double]"0.30470667752978903"
$ws.Range("AS63").Value = [double]"0.34105242119623103"
$ws.Range("AT63").Value = [double]"0.19697008382799999"
$ws.Range("AU63").Value = [double]"0.13069833122553401"
$ws.Range("AV63").Value = "across_condition"
$ws.Range("AR64").Value = [double]"0.49386386018753398"
$ws.Range("AS64").Value = [double]"0.39485681153025598"
$ws.Range("AT64").Value = [double]"0.29031800265242202"
$ws.Range("AU64").Value = [double]"0.231727445795535"
$ws.Range("AV64").Value = "across_condition"
$ws.Range("AR65").Value = [double]"0.20668499018887199"
$ws.Range("AS65").Value = [double]"0.22906731319409901"
$ws.Range("AT65").Value = [double]"0.27499469118066"
$ws.Range("AU65").Value = [double]"0.14656653254521099"
$ws.Range("AV65").Value = "across_condition"
$ws.Range("AR66").Value = [double]"0.106869146175161"
$ws.Range("AS66").Value = [double]"9.2900217465541904E-2"
$ws.Range("AT66").Value = [double]"0.12820125328458901"
$ws.Range("AU66").Value = [double]"6.0516094803882002E-2"
$ws.Range("AV66").Value = "across_condition"
$ws.Range("AR67").Value = [double]"0.48498985972303099"
$ws.Range("AS67").Value = [double]"0.50273515618515596"
$ws.Range("AT67").Value = [double]"0.47525301388618602"
$ws.Range("AU67").Value = [double]"0.30802403317891403"
$ws.Range("AV67").Value = "across_condition"
$ws.Range("AR68").Value = [double]"0.48498985972303099"
$ws.Range("AS68").Value = [double]"0.50273515618515596"
$ws.Range("AT68").Value = [double]"0.47525301388618602"
$ws.Range("AU68").Value = [double]"0.30802403317891403"
$ws.Range("AV68").Value = "across_condition"
$ws.Range("AR69").Value = [double]"7.7986698352480399E-2"
$ws.Range("AS69").Value = [double]"0.21351360929719199"
$ws.Range("AT69").Value = [double]"3.9737395730011699E-2"
$ws.Range("AU69").Value = [double]"4.6662803949643797E-2"
$ws.Range("AV69").Value = "Stim8hr"
$ws.Range("AR70").Value = [double]"0.20085245145347499"
$ws.Range("AS70").Value = [double]"0.325302403689186"
$ws.Range("AT70").Value = [double]"0.13201608881397101"
$ws.Range("AU70").Value = [double]"0.111923988620486"
$ws.Range("AV70").Value = "across_condition"
$ws.Range("AR71").Value = [double]"7.5905858579455199E-2"
$ws.Range("AS71").Value = [double]"6.2400265196790899E-2"
$ws.Range("AT71").Value = [double]"3.7098317073143101E-2"
$ws.Range("AU71").Value = [double]"3.5659291108160997E-2"
$ws.Range("AV71").Value = "across_condition"
$ws.Range("AR72").Value = [double]"0.17148867371832199"
$ws.Range("AS72").Value = [double]"0.23497401176868499"
$ws.Range("AT72").Value = [double]"0.151214584771793"
$ws.Range("AU72").Value = [double]"0.13942275295025799"
$ws.Range("AV72").Value = "across_condition"
$ws.Range("AR73").Value = [double]"9.2598439582371794E-2"
$ws.Range("AS73").Value = [double]"0.14195976445227501"
$ws.Range("AT73").Value = [double]"0.13002722105876399"
$ws.Range("AU73").Value = [double]"0.14542747893691599"
$ws.Range("AV73").Value = "across_condition"
$ws.Range("AR74").Value = [double]"0.24820649282461499"
$ws.Range("AS74").Value = [double]"7.50502611665297E-2"
$ws.Range("AT74").Value = [double]"3.4198618483775499E-2"
$ws.Range("AU74").Value = [double]"3.5390860630754202E-2"
$ws.Range("AV74").Value = "Rest"
$ws.Range("AR75").Value = [double]"4.14951162538763E-2"
$ws.Range("AS75").Value = [double]"7.3945953394869396E-2"
$ws.Range("AT75").Value = [double]"0.103354162453844"
$ws.Range("AU75").Value = [double]"5.3846482451243599E-2"
$ws.Range("AV75").Value = "across_condition"
$ws.Range("AR76").Value = [double]"0.225857347223547"
$ws.Range("AS76").Value = [double]"0.319560876850779"
$ws.Range("AT76").Value = [double]"0.286553898566551"
$ws.Range("AU76").Value = [double]"0.112902830722226"
$ws.Range("AV76").Value = "across_condition"
$ws.Range("AR77").Value = [double]"0.33022771554444202"
$ws.Range("AS77").Value = [double]"0.12308387274312101"
$ws.Range("AT77").Value = [double]"5.7215501257170898E-2"
$ws.Range("AU77").Value = [double]"7.4355333464021095E-2"
$ws.Range("AV77").Value = "Rest"
$ws.Range("AR78").Value = [double]"3.3634521580691697E-2"
$ws.Range("AS78").Value = [double]"3.7936240309417099E-2"
$ws.Range("AT78").Value = [double]"0.15485081356272401"
$ws.Range("AU78").Value = [double]"2.3850552053401201E-2"
$ws.Range("AV78").Value = "Stim48hr"
$ws.Range("AR79").Value = [double]"4.4033365616632902E-2"
$ws.Range("AS79").Value = [double]"3.5389745744152698E-2"
$ws.Range("AT79").Value = [double]"0.25573679194922999"
$ws.Range("AU79").Value = [double]"4.1339331645289197E-2"
$ws.Range("AV79").Value = "Stim48hr"
$ws.Range("AR80").Value = [double]"0.15227071501335401"
$ws.Range("AS80").Value = [double]"0.17200271580357601"
$ws.Range("AT80").Value = [double]"0.27844877914084598"
$ws.Range("AU80").Value = [double]"0.121569601813637"
$ws.Range("AV80").Value = "across_condition"
$ws.Range("AR81").Value = [double]"0.30491986556016298"
$ws.Range("AS81").Value = [double]"0.21322057865573801"
$ws.Range("AT81").Value = [double]"0.184010658074235"
$ws.Range("AU81").Value = [double]"9.9303900106460397E-2"
$ws.Range("AV81").Value = "across_condition"
$ws.Range("AR82").Value = [double]"0.147428943789184"
$ws.Range("AS82").Value = [double]"0.13566781398655101"
$ws.Range("AT82").Value = [double]"4.3363498081228698E-2"
$ws.Range("AU82").Value = [double]"4.2872357048438899E-2"
$ws.Range("AV82").Value = "Rest_Stim8hr"
$ws.Range("AR83").Value = [double]"0.351597677908358"
$ws.Range("AS83").Value = [double]"0.37319360583713501"
$ws.Range("AT83").Value = [double]"0.378375731639229"
$ws.Range("AU83").Value = [double]"0.27571699939382699"
$ws.Range("AV83").Value = "across_condition"
$ws.Range("AR84").Value = [double]"3.4228870007714102E-2"
$ws.Range("AS84").Value = [double]"5.6747444631806698E-2"
$ws.Range("AT84").Value = [double]"0.167098248494992"
$ws.Range("AU84").Value = [double]"5.1071616492417801E-2"
$ws.Range("AV84").Value = "Stim48hr"
$ws.Range("AR85").Value = [double]"0.29580802265957301"
$ws.Range("AS85").Value = [double]"0.28576289465232102"
$ws.Range("AT85").Value = [double]"0.322213258458517"
$ws.Range("AU85").Value = [double]"0.18048573635029"
$ws.Range("AV85").Value = "across_condition"
$ws.Range("AR86").Value = [double]"0.38221840022190601"
$ws.Range("AS86").Value = [double]"0.17915470912607501"
$ws.Range("AT86").Value = [double]"5.5105186947498601E-2"
$ws.Range("AU86").Value = [double]"9.7535399270617107E-2"
$ws.Range("AV86").Value = "Rest"
$ws.Range("AR87").Value = [double]"1.9649069258256501E-2"
$ws.Range("AS87").Value = [double]"7.9786296373786994E-2"
$ws.Range("AT87").Value = [double]"0.217369418152452"
$ws.Range("AU87").Value = [double]"3.2432790130622703E-2"
$ws.Range("AV87").Value = "Stim48hr"
$ws.Range("AR88").Value = [double]"3.3535806629703298E-2"
$ws.Range("AS88").Value = [double]"0.132726455736841"
$ws.Range("AT88").Value = [double]"5.7284502761605599E-2"
$ws.Range("AU88").Value = [double]"1.1707811967664201E-2"
$ws.Range("AV88").Value = "Stim8hr"
$ws.Range("AR89").Value = [double]"0.24763360996150599"
$ws.Range("AS89").Value = [double]"0.16967093068265501"
$ws.Range("AT89").Value = [double]"9.6916792077621106E-2"
$ws.Range("AU89").Value = [double]"0.10598591727578301"
$ws.Range("AV89").Value = "across_condition"
$ws.Range("AR90").Value = [double]"0.166182611605936"
$ws.Range("AS90").Value = [double]"0.18139967642655899"
$ws.Range("AT90").Value = [double]"4.6886502848284398E-2"
$ws.Range("AU90").Value = [double]"6.3457907626463E-2"
$ws.Range("AV90").Value = "Rest_Stim8hr"
$ws.Range("AR91").Value = [double]"0.262674573072117"
$ws.Range("AS91").Value = [double]"6.4061365573942194E-2"
$ws.Range("AT91").Value = [double]"1.0193169502506099E-2"
$ws.Range("AU91").Value = [double]"7.34077974354818E-2"
$ws.Range("AV91").Value = "Rest"
$ws.Range("AR92").Value = [double]"1.48081695720771E-2"
$ws.Range("AS92").Value = [double]"4.4723679692697797E-2"
$ws.Range("AT92").Value = [double]"0.12727436373547099"
$ws.Range("AU92").Value = [double]"9.0975251131272996E-3"
$ws.Range("AV92").Value = "Stim48hr"
$ws.Range("AR93").Value = [double]"0.36710569412349803"
$ws.Range("AS93").Value = [double]"0.234583503166892"
$ws.Range("AT93").Value = [double]"0.26253658896203302"
$ws.Range("AU93").Value = [double]"0.16575579696873999"
$ws.Range("AV93").Value = "across_condition"
$ws.Range("AR94").Value = [double]"0.34843759677928199"
$ws.Range("AS94").Value = [double]"0.13670138676129701"
$ws.Range("AT94").Value = [double]"3.72918494079734E-2"
$ws.Range("AU94").Value = [double]"0.102169260643659"
$ws.Range("AV94").Value = "Rest"
$ws.Range("AR95").Value = [double]"5.153647668585E-2"
$ws.Range("AS95").Value = [double]"0.48299508452469297"
$ws.Range("AT95").Value = [double]"0.25424943538616701"
$ws.Range("AU95").Value = [double]"0.13338368097925801"
$ws.Range("AV95").Value = "Stim8hr_Stim48hr"
$ws.Range("AR96").Value = [double]"6.9025630436642393E-2"
$ws.Range("AS96").Value = [double]"5.62469308217197E-2"
$ws.Range("AT96").Value = [double]"0.18120608337073599"
$ws.Range("AU96").Value = [double]"5.5648682688465803E-2"
$ws.Range("AV96").Value = "Stim48hr"
$ws.Range("AR97").Value = [double]"3.4628360437118401E-2"
$ws.Range("AS97").Value = [double]"5.1071018977830701E-2"
$ws.Range("AT97").Value = [double]"0.227014918883545"
$ws.Range("AU97").Value = [double]"2.93313319513713E-2"
$ws.Range("AV97").Value = "Stim48hr"
$ws.Range("AR98").Value = [double]"0.11623028447802999"
$ws.Range("AS98").Value = [double]"0.19251221958172901"
$ws.Range("AT98").Value = [double]"0.101641784321131"
$ws.Range("AU98").Value = [double]"8.3587140933692697E-2"
$ws.Range("AV98").Value = "across_condition"
$ws.Range("AR99").Value = [double]"0.10389471323231"
$ws.Range("AS99").Value = [double]"0.10780767189316"
$ws.Range("AT99").Value = [double]"5.8396497820150897E-2"
$ws.Range("AU99").Value = [double]"4.0871517004645398E-2"
$ws.Range("AV99").Value = "across_condition"
$ws.Range("AR100").Value = [double]"0.168417378499611"
$ws.Range("AS100").Value = [double]"0.220596822680375"
$ws.Range("AT100").Value = [double]"0.34881713830814798"
$ws.Range("AU100").Value = [double]"0.119637691225501"
$ws.Range("AV100").Value = "across_condition"
$ws.Range("AR101").Value = [double]"6.8010509514295797E-2"
$ws.Range("AS101").Value = [double]"9.7937928666221197E-2"
$ws.Range("AT101").Value = [double]"0.27616319735343198"
$ws.Range("AU101").Value = [double]"8.3176783990759895E-2"
$ws.Range("AV101").Value = "Stim48hr"
$ws.Range("AR102").Value = [double]"4.4663060017211903E-2"
$ws.Range("AS102").Value = [double]"0.13069514557324199"
$ws.Range("AT102").Value = [double]"0.192627226148222"
$ws.Range("AU102").Value = [double]"4.6344766647650103E-2"
$ws.Range("AV102").Value = "Stim48hr"
$ws.Range("AR103").Value = [double]"0.36163734582123103"
$ws.Range("AS103").Value = [double]"0.14071592248065701"
$ws.Range("AT103").Value = [double]"5.1399245015027101E-2"
$ws.Range("AU103").Value = [double]"9.0572324812866298E-2"
$ws.Range("AV103").Value = "Rest"
$ws.Range("AR104").Value = [double]"0.111094181769095"
$ws.Range("AS104").Value = [double]"8.5692970171344696E-2"
$ws.Range("AT104").Value = [double]"0.37466147687576301"
$ws.Range("AU104").Value = [double]"0.10781378316586"
$ws.Range("AV104").Value = "Stim48hr"
$ws.Range("AR105").Value = [double]"-1.17711935912774E-2"
$ws.Range("AS105").Value = [double]"3.8051300928665E-3"
$ws.Range("AT105").Value = [double]"0.42137823431745702"
$ws.Range("AU105").Value = [double]"-2.6217859603379001E-3"
$ws.Range("AV105").Value = "Stim48hr"
$ws.Range("AR106").Value = [double]"0.48414901971790902"
$ws.Range("AS106").Value = [double]"0.18668107351959901"
$ws.Range("AT106").Value = [double]"3.7648907670759797E-2"
$ws.Range("AU106").Value = [double]"7.2489908369905307E-2"
$ws.Range("AV106").Value = "Rest"
$ws.Range("AR107").Value = [double]"0.36295124321286998"
$ws.Range("AS107").Value = [double]"9.8820120714134305E-2"
$ws.Range("AT107").Value = [double]"2.1683013812557401E-2"
$ws.Range("AU107").Value = [double]"5.8439687501210598E-2"
$ws.Range("AV107").Value = "Rest"
$ws.Range("AR108").Value = [double]"0.28307013297014499"
$ws.Range("AS108").Value = [double]"7.9958533610863206E-2"
$ws.Range("AT108").Value = [double]"6.0316869610292E-2"
$ws.Range("AU108").Value = [double]"6.9853487169400097E-2"
$ws.Range("AV108").Value = "Rest"
$ws.Range("AR109").Value = [double]"7.4705122379974703E-2"
$ws.Range("AS109").Value = [double]"7.0581344667149298E-2"
$ws.Range("AT109").Value = [double]"0.35315478797448602"
$ws.Range("AU109").Value = [double]"0.100313280243879"
$ws.Range("AV109").Value = "Stim48hr"
$ws.Range("AR110").Value = [double]"0.34363766792706701"
$ws.Range("AS110").Value = [double]"0.24107893908767"
$ws.Range("AT110").Value = [double]"0.38816597650156698"
$ws.Range("AU110").Value = [double]"0.15852204065464601"
$ws.Range("AV110").Value = "across_condition"
$ws.Range("AR111").Value = [double]"3.9248112075492801E-2"
$ws.Range("AS111").Value = [double]"5.3211923173604299E-2"
$ws.Range("AT111").Value = [double]"0.29121290438076097"
$ws.Range("AU111").Value = [double]"5.2853807980403497E-2"
$ws.Range("AV111").Value = "Stim48hr"
$ws.Range("AR112").Value = [double]"5.73771025867248E-2"
$ws.Range("AS112").Value = [double]"5.5012014814709298E-2"
$ws.Range("AT112").Value = [double]"0.25209150350947102"
$ws.Range("AU112").Value = [double]"5.7645210367781502E-2"
$ws.Range("AV112").Value = "Stim48hr"
$ws.Range("AR113").Value = [double]"0.18500936463165399"
$ws.Range("AS113").Value = [double]"0.156122485146936"
$ws.Range("AT113").Value = [double]"0.381914944126369"
$ws.Range("AU113").Value = [double]"8.8984295372876193E-2"
$ws.Range("AV113").Value = "Stim48hr"

# --- View state: the user scrolled/selected the new columns after adding them ---
$ws.Range("AR1:AV1048576").Select()
